$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A40").Value = "간격"
$ws.Range("B40").Value = "간격,당첨번호 간의 번호 간격을 계산하고 평균을 분석합니다."

$ws.Range("F41").Select()
